# repull data, push all data, mean calculation
# Update the dSF (column F) values to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    4  = 8
    8  = 2
    9  = -9
    10 = 1
    11 = 0
    19 = -7
    20 = -1
    22 = 0
    23 = -1
    24 = -3
    25 = -2
    28 = 1
    29 = -2
    30 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
